# Insert a new column at the beginning, shifting existing data (A:E) to (B:F),
# then populate the new column A with an "ID" header and per-row labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A:E shift to B:F.
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold/border/centered) from B1 into the new A1 cell.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "ID" column with header + row labels.
$ids = @("ID","Hb 2","Hb 3","S 24","S 28","Hb 107","Hb 66","Hb 69","Hb 95","Hb 99","Hb 92","Hb 40","Hb 41","S 11","Hb 57","S 21","S 22","S 3","S 4","S 5","Hb 74","Hb 79","Hb 32","S 15","S 16")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
